# Rename sheet "August" -> "November" and update the Print_Area defined name accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "November"
$ws.PageSetup.PrintArea = "`$A`$1:`$J`$25"

# Row 10: add a "26" value for the second half of the month in column F.
$ws.Range("F10").Value = 26

# Row 11: time-slot text. E11 switches from the 1pm-4pm slot to a new 12:30pm-6:30pm slot,
# and F11 picks up the same new slot text (previously blank).
$ws.Range("E11").Value = "12:30 PM To 06:30 PM"
$ws.Range("F11").Value = "12:30 PM To 06:30 PM"

# Row 12: add classes count for column F; I12 total (SUM(D12:H12)) recalculates automatically.
$ws.Range("F12").Value = 6

# Rows 14-17: column F picks up explicit 0 values (previously blank).
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 0

# Row 18: E18 goes from 3 to 6, F18 goes from blank to 6.
$ws.Range("E18").Value = 6
$ws.Range("F18").Value = 6

# Row 19: E19 goes from 3 to 6, F19 goes from blank to 6.
$ws.Range("E19").Value = 6
$ws.Range("F19").Value = 6

# Row 20: E20 goes from 3 to 6, F20 goes from blank to 6.
$ws.Range("E20").Value = 6
$ws.Range("F20").Value = 6

# Row 21: E21 goes from 3 to 6, F21 goes from blank to 0.
$ws.Range("E21").Value = 6
$ws.Range("F21").Value = 0

# Sheet view: scroll/selection moved from A7/G19 to A8/I19.
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("I19").Select() | Out-Null
